$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G4").Value = "2016-08-29 02:45:19"

$wsZhCn.Range("H4").Value = "2016-08-29 02:45:14"
$wsZhCn.Range("K4").Value = "2016-08-29 02:45:32"

$wsDeDe.Range("K4").Value = "2016-08-29 02:45:39"
